$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (N_Calib_1=20, N_Calib_2=40)
$ws.Range("C2").Value = -0.3309010421165708
$ws.Range("D2").Value = 0.7438502100025386

# Row 3 (N_Calib_1=20, N_Calib_2=60)
$ws.Range("C3").Value = 0.738802337744134
$ws.Range("D3").Value = 0.4678365618825044

# Row 4 (N_Calib_1=20, N_Calib_2=100)
$ws.Range("C4").Value = 2.467652683724732
$ws.Range("D4").Value = 0.02185422460194419

# Row 5 (N_Calib_1=20, N_Calib_2=200)
$ws.Range("C5").Value = 3.840728866733985
$ws.Range("D5").Value = 0.000888710136857318

# Row 6 (N_Calib_1=40, N_Calib_2=60)
$ws.Range("C6").Value = 0.7280352027440316
$ws.Range("D6").Value = 0.4742702704389119

# Row 7 (N_Calib_1=40, N_Calib_2=100)
$ws.Range("C7").Value = 3.002372990719762
$ws.Range("D7").Value = 0.006558448967928898

# Row 8 (N_Calib_1=40, N_Calib_2=200)
$ws.Range("C8").Value = 3.191883911006894
$ws.Range("D8").Value = 0.004211822780769703

# Row 9 (N_Calib_1=60, N_Calib_2=100)
$ws.Range("C9").Value = 1.742158862356944
$ws.Range("D9").Value = 0.09544562098095866
$ws.Range("G9").Value = "No"

# Row 10 (N_Calib_1=60, N_Calib_2=200)
$ws.Range("C10").Value = 2.853251771445125
$ws.Range("D10").Value = 0.00924370360146165

# Row 11 (N_Calib_1=100, N_Calib_2=200)
$ws.Range("C11").Value = -0.1504460470507195
$ws.Range("D11").Value = 0.8817835737091928
